$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue 'D2' '30.791.15'
Set-TextValue 'E2' '  +0.53%  '
Set-TextValue 'D3' '1.889.30'
Set-TextValue 'E3' '  +1.38%  '
Set-TextValue 'E4' '  +0.26%  '
Set-TextValue 'D5' '239.98'
Set-TextValue 'E5' '  +2.32%  '
Set-TextValue 'E6' '  +0.27%  '
Set-TextValue 'D7' '0.4784'
Set-TextValue 'E7' '  +1.98%  '
Set-TextValue 'D8' '0.2974'
Set-TextValue 'E8' '  +7.85%  '
Set-TextValue 'D9' '0.06660'
Set-TextValue 'E9' '  +4.74%  '
Set-TextValue 'D10' '18.65'
Set-TextValue 'E10' '  +6.75%  '
Set-TextValue 'D11' '100.36'
Set-TextValue 'E11' '  +18.07%  '
Set-TextValue 'D12' '1.879.03'
Set-TextValue 'E12' '  +0.84%  '
Set-TextValue 'D13' '0.07560'
Set-TextValue 'E13' '  +1.32%  '
Set-TextValue 'D14' '5.153'
Set-TextValue 'E14' '  +3.89%  '
Set-TextValue 'D15' '0.6613'
Set-TextValue 'E15' '  +4.82%  '
Set-TextValue 'D16' '303.39'
Set-TextValue 'E16' '  +25.24%  '
Set-TextValue 'D17' '30.771.76'
Set-TextValue 'E17' '  +0.69%  '
Set-TextValue 'E18' '  +3.49%  '
Set-TextValue 'D19' '1.001'
Set-TextValue 'E19' '  +0.28%  '
Set-TextValue 'D20' '0.000007615'
Set-TextValue 'E20' '  +3.41%  '
Set-TextValue 'D21' '2.124.04'
Set-TextValue 'E21' '  +1.51%  '
Set-TextValue 'E22' '  +0.14%  '
Set-TextValue 'D23' '5.161'
Set-TextValue 'E23' '  +3.51%  '
Set-TextValue 'D24' '6.217'
Set-TextValue 'E24' '  +4.25%  '
Set-TextValue 'D25' '9.319'
Set-TextValue 'E25' '  +0.57%  '
Set-TextValue 'D26' '167.96'
Set-TextValue 'E26' '  +0.65%  '
Set-TextValue 'D27' '20.43'
Set-TextValue 'E27' '  +12.30%  '
Set-TextValue 'D28' '1.953'
Set-TextValue 'E28' '  +3.48%  '
Set-TextValue 'D29' '0.1130'
Set-TextValue 'E29' '  +10.02%  '
Set-TextValue 'E30' '  -1.56%  '
Set-TextValue 'E31' '  +1.56%  '
Set-TextValue 'D32' '3.996'
Set-TextValue 'D33' '0.05096'
Set-TextValue 'E33' '  +3.35%  '
Set-TextValue 'B34' 'ImmutableX'
Set-TextValue 'C34' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D34' '0.7534'
Set-TextValue 'E34' '  +6.02%  '
Set-TextValue 'B35' 'ARBITRUM'
Set-TextValue 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '1.164'
Set-TextValue 'E35' '  +1.12%  '
Set-TextValue 'E36' '  +0.76%  '
Set-TextValue 'D37' '0.01981'
Set-TextValue 'E37' '  +3.56%  '
Set-TextValue 'D38' '2.710'
Set-TextValue 'E38' '  +0.63%  '
Set-TextValue 'D39' '2.064'
Set-TextValue 'E39' '  +4.12%  '
Set-TextValue 'D40' '0.8966'
Set-TextValue 'E40' '  +1.54%  '
Set-TextValue 'D41' '108.19'
Set-TextValue 'E41' '  +2.23%  '
Set-TextValue 'E42' '  +0.29%  '
Set-TextValue 'D43' '0.4201'
Set-TextValue 'E43' '  +2.49%  '
Set-TextValue 'D44' '5.650'
Set-TextValue 'E44' '  +1.88%  '
Set-TextValue 'D45' '67.67'
Set-TextValue 'E45' '  +9.19%  '
Set-TextValue 'D46' '7.382'
Set-TextValue 'E46' '  +2.52%  '
Set-TextValue 'D47' '9.071'
Set-TextValue 'E47' '  +5.40%  '
Set-TextValue 'D48' '0.1233'
Set-TextValue 'E48' '  -0.42%  '
Set-TextValue 'D49' '34.98'
Set-TextValue 'E49' '  +3.75%  '
Set-TextValue 'E50' '  +1.60%  '
Set-TextValue 'E51' '  +1.74%  '
